# Fix formatting bugs introduced when scraping floating point numbers /
# proveedor names from the source system.
#
# 1) A handful of "Razon social" / "Nombre Fantasia" entries used a comma as
#    a separator between multiple people/partners in a company name
#    (e.g. "PARPAGNOLI, PEDRO RICARDO"). Those commas are replaced with
#    periods, and a couple of "S.H." suffixes are normalised to "SH".
#
# 2) Every value in the "Importe" column (H) was scraped as Spanish/
#    Argentina-formatted text (thousands separator "." and decimal
#    separator ","), e.g. "13.680,00". The fix rewrites these as plain
#    decimal text with a "." decimal separator and no thousands grouping,
#    e.g. "13680.00" -- while keeping the cell content as literal text
#    (Excel would otherwise silently re-interpret a clean numeric string
#    typed into a General-formatted cell as a real number, truncating
#    trailing zeros such as "13680.00" -> 13680). Prefixing the new value
#    with a leading apostrophe forces Excel to store it verbatim as text,
#    matching the original "text-only" shared-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Proveedor / company name punctuation fixes ---------------------
$nameFixes = @{
    "E103" = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
    "E111" = "PARPAGNOLI. PEDRO RICARDO"
    "F111" = "PARPAGNOLI. PEDRO RICARDO ( SP )"
    "E145" = "GIMENEZ ANIBAL. FALISTOCCO MARISA DANIELA SH"
    "E184" = "PARPAGNOLI. PEDRO RICARDO"
    "F184" = "PARPAGNOLI. PEDRO RICARDO ( SP )"
    "E200" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E216" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
}

foreach ($addr in $nameFixes.Keys) {
    $ws.Range($addr).Value = $nameFixes[$addr]
}

# --- 2) "Importe" column (H) number formatting fix ----------------------
# Data rows are 2..286 ("H1" is the "Importe" header).
$firstRow = 2
$lastRow = 286
$col = 8  # column H

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $old = $cell.Text
    $new = $old.Replace(".", "").Replace(",", ".")
    # Leading apostrophe => keep as text instead of letting Excel coerce
    # the clean numeric-looking string into a real Number value.
    $cell.Value = "'" + $new
}
